# Insert a new weekly record at row 396 ("Hortaliza, Vega Central Mapocho de
# Santiago - Haba"), pushing the existing rows 396-429 down to 397-430.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(396).Insert()

$ws.Cells.Item(396, 1).Value  = 9
$ws.Cells.Item(396, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(396, 3).Value  = "Metropolitana"
$ws.Cells.Item(396, 4).Value  = 45265
$ws.Cells.Item(396, 5).Value  = 13
$ws.Cells.Item(396, 6).Value  = 100112026
$ws.Cells.Item(396, 7).Value  = "Haba"
$ws.Cells.Item(396, 8).Value  = "Sin especificar"
$ws.Cells.Item(396, 9).Value  = "Primera"
$ws.Cells.Item(396, 10).Value = 70
$ws.Cells.Item(396, 11).Value = 11000
$ws.Cells.Item(396, 12).Value = 13000
$ws.Cells.Item(396, 13).Value = 12000
$ws.Cells.Item(396, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(396, 15).Value = "Región del Maule"
$ws.Cells.Item(396, 16).Value = 480
$ws.Cells.Item(396, 17).Value = 25
$ws.Cells.Item(396, 18).Value = "Hortaliza"
